$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 30: new literature entry (Hilbert/Redmiles) ---------------------
# Columns A, B, C, F, G, H simply inherit the sheet's per-column default
# style (text / text / number / number / text / text) the moment a value is
# written into a previously-empty cell, so a plain value assignment is
# enough for those. D needs the explicit short-date format, and E already
# carries the wrap-text style from the template fill-down, so only its
# value needs to be set.
$ws.Range("A30").Value2 = "Separating the wheat from the chaff in Internet-mediated user feedback expectation-driven event monitoring"
$ws.Range("B30").Value2 = "David M. Hilbert and David F. Redmiles"
$ws.Range("C30").Value2 = 1999
$ws.Range("D30").Value2 = 44113
$ws.Range("D30").NumberFormat = "m/d/yy"
$ws.Range("E30").Value2 = "Beschreibt welche Arten es für Informationserhebung gibt und wie man an die relevanten Daten kommt"
$ws.Range("F30").Value2 = 5
$ws.Range("G30").Value2 = "TU-Bib runterladen"
$ws.Range("H30").Value2 = "http://citeseerx.ist.psu.edu/viewdoc/download;jsessionid=6FA1A0A2B3361C59E5220375CD895D31?doi=10.1.1.51.9375&rep=rep1&type=pdf"
$ws.Rows(30).RowHeight = 43.2

# --- Row 31: new literature entry (Wobser) --------------------------------
$ws.Range("A31").Value2 = "Produktentwicklung in Kooperation mit Anwendern"
$ws.Range("B31").Value2 = "Gunther Wobser"
$ws.Range("C31").Value2 = 37687
$ws.Range("C31").NumberFormat = "m/d/yy"
$ws.Range("D31").Value2 = 44113
$ws.Range("D31").NumberFormat = "m/d/yy"
$ws.Range("F31").Value2 = 3
$ws.Range("G31").Value2 = "FH-VPN ?"
$ws.Range("H31").Value2 = "https://link.springer.com/book/10.1007/978-3-322-81517-0"

# --- Selection, as left by the author after editing -----------------------
$ws.Range("C31").Select() | Out-Null
